# Generate Report for Archive
# - Updates the localization status from "Ready for handoff" to "In Translation"
#   on every sheet that references it (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - Re-fits the now-narrower status columns to their new (shorter) content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status value everywhere it appears.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# The status columns were sized for the old (longer) text; shrink them to
# fit the new, shorter text now that the report has been regenerated.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
